$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated GDP figures for existing rows (B12:B30)
$updatedValues = @{
    12 = 208125.35500000001
    13 = 214103.45600000001
    14 = 224545.62899999999
    15 = 239911.62400000001
    16 = 257637.92199999999
    17 = 269153.32
    18 = 277955.08100000001
    19 = 273531.55099999998
    20 = 266522.79700000002
    21 = 272123.15899999999
    22 = 281302.90500000003
    23 = 292330.81300000002
    24 = 305248.97499999998
    25 = 325002.75699999998
    26 = 347744.95699999999
    27 = 369995.897
    28 = 389471.81300000002
    29 = 411660.04499999998
    30 = 438936.81400000001
}

foreach ($row in $updatedValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $updatedValues[$row]
}

# New row 31: 2020-01-01 observation date (serial 43831) and GDP value
$ws.Cells.Item(31, 1).Value = 43831
$ws.Cells.Item(31, 1).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(31, 2).Value = 425443.80699999997
$ws.Cells.Item(31, 2).NumberFormat = "0.000"

# Match the author's updated selection (full column A:B select, active cell A12)
$ws.Cells.Item(12, 1).Select()
$ws.Range("A1:B1048576").Select()
